# Weekly update: add two new daily price rows (Fruta / hortaliza, semanal)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at sheet row 3 (2023-04-19) ---------------------------
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 45035
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100112039
$ws.Cells.Item(3, 7).Value = "Ciboulette"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 1100
$ws.Cells.Item(3, 11).Value = 2000
$ws.Cells.Item(3, 12).Value = 2500
$ws.Cells.Item(3, 13).Value = 2250
$ws.Cells.Item(3, 14).Value = "`$/docena de atados"
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 750
$ws.Cells.Item(3, 17).Value = 3
$ws.Cells.Item(3, 18).Value = "Hortaliza"

# --- Insert new row at sheet row 20 (2023-04-18) ---------------------------
$ws.Rows.Item(20).Insert()
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 45034
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 100112039
$ws.Cells.Item(20, 7).Value = "Ciboulette"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 1100
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 2500
$ws.Cells.Item(20, 13).Value = 2250
$ws.Cells.Item(20, 14).Value = "`$/docena de atados"
$ws.Cells.Item(20, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(20, 16).Value = 750
$ws.Cells.Item(20, 17).Value = 3
$ws.Cells.Item(20, 18).Value = "Hortaliza"
